# Apply cryptocurrency price/volume update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "49.232.40"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.89%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.629.17"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.45%  "

$ws.Range("E4").Value = "  +0.07%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "111.99"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.52%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "322.91"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.33%  "

$ws.Range("E7").Value = "  -1.30%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  -3.01%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "39.73"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.31%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "19.73"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -5.07%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0810"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.03%  "

$ws.Range("E13").Value = "  +1.08%  "

$ws.Range("E14").Value = "  -0.11%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.038.89"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.37%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "2.637.99"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.41%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.858"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.67%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "49.176.92"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.91%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.01"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.58%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.92"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -3.16%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.69"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.81%  "

$ws.Range("E22").Value = "  -0.60%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "269.77"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -3.82%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "68.54"
$c.Style = "Normal"

$ws.Range("E25").Value = "  -1.73%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "26.19"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.48%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.08%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "10.22"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +2.50%  "

$ws.Range("E29").Value = "  -0.29%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "35.10"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.91%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.138"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -4.06%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "49.51"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.57%  "

$ws.Range("E33").Value = "  +0.44%  "

$ws.Range("E34").Value = "  +1.18%  "

$ws.Range("E35").Value = "  -0.29%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "19.05"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -3.47%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "4.95"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +4.37%  "

$ws.Range("E38").Value = "  -0.67%  "

$ws.Range("E39").Value = "  +1.10%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "126.04"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.82%  "

$ws.Range("E41").Value = "  -1.65%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "22.18"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.28%  "

$ws.Range("E43").Value = "  -4.26%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0318"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.86%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.065.97"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.77%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.16"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +7.05%  "

$ws.Range("E47").Value = "  -4.16%  "

$ws.Range("E48").Value = "  -5.00%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "8.90"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.53%  "

$ws.Range("E50").Value = "  -2.79%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "58.85"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.91%  "
